$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update "想去人数" (F column) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 542
$ws1.Range("F5").Value = 523
$ws1.Range("F6").Value = 292
$ws1.Range("F7").Value = 2636
$ws1.Range("F9").Value = 7248
$ws1.Range("F11").Value = 452
$ws1.Range("F13").Value = 171

# Sheet "全部类型" (sheet4) - update "想去人数" (F column) counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 542
$ws4.Range("F5").Value = 523
$ws4.Range("F6").Value = 292
$ws4.Range("F9").Value = 2636
$ws4.Range("F11").Value = 7248
$ws4.Range("F13").Value = 452
$ws4.Range("F17").Value = 171
